# Remove dependency on matsim run output config file:
# add a new "networkCrs" parameter row to the scenario_info sheet
# (global group), since the network CRS can no longer be read from
# the matsim run output config and must now be supplied directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# Insert a new blank row above the current row 6 ("sampleSize"),
# shifting it (and everything below) down by one.
$ws.Rows("6:6").Insert()

# Populate the new row with the networkCrs parameter.
$ws.Range("A6").Value2 = "global"
$ws.Range("B6").Value2 = "networkCrs"
$ws.Range("C6").Value2 = "EPSG:25832"
$ws.Range("E6").Value2 = "The coordinate reference system of the network"

# The autofilter range needs to grow to cover the extra row. Turn it
# off and reapply it over the full, now-larger, table range.
$ws.AutoFilterMode = $false
$ws.Range("A1:E25").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new
# autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "scenario_info!_FilterDatabase") {
        $n.RefersTo = "=scenario_info!`$A`$1:`$E`$25"
    }
}

# Reflect the last user selection recorded for this sheet.
$ws.Activate()
$ws.Range("B11").Select()
